$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Updated stats for rows 2-11 (NoOfSessions, PointsBonusTotal, AvgSessPoints)
$ws.Range("D2").Value = 127
$ws.Range("I2").Value = 460
$ws.Range("J2").Value = 3.62

$ws.Range("D3").Value = 146
$ws.Range("I3").Value = 525
$ws.Range("J3").Value = 3.6

$ws.Range("D4").Value = 196
$ws.Range("I4").Value = 733
$ws.Range("J4").Value = 3.74

$ws.Range("D5").Value = 108
$ws.Range("I5").Value = 394
$ws.Range("J5").Value = 3.65

$ws.Range("D6").Value = 43
$ws.Range("I6").Value = 158
$ws.Range("J6").Value = 3.67

$ws.Range("D7").Value = 147
$ws.Range("I7").Value = 626
$ws.Range("J7").Value = 4.26

$ws.Range("D8").Value = 201
$ws.Range("I8").Value = 705
$ws.Range("J8").Value = 3.51

$ws.Range("D9").Value = 87
$ws.Range("J9").Value = 3.38

$ws.Range("D10").Value = 204
$ws.Range("I10").Value = 839
$ws.Range("J10").Value = 4.11

$ws.Range("D11").Value = 124
$ws.Range("I11").Value = 487
